$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 682
$ws.Range("F6").Value = 610
$ws.Range("F9").Value = 38
$ws.Range("F10").Value = 6093
$ws.Range("F11").Value = 664
$ws.Range("F12").Value = 1066
$ws.Range("F13").Value = 27
$ws.Range("F14").Value = 264
$ws.Range("F17").Value = 601
$ws.Range("F18").Value = 1034
$ws.Range("F19").Value = 56
$ws.Range("F20").Value = 38
$ws.Range("F21").Value = 218
$ws.Range("F22").Value = 1363
$ws.Range("F24").Value = 1030
$ws.Range("F25").Value = 75
$ws.Range("F26").Value = 2095
$ws.Range("F27").Value = 204
$ws.Range("F28").Value = 41
$ws.Range("F31").Value = 3385

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 349
$ws.Range("G6").Value = 180
$ws.Range("F10").Value = 665
$ws.Range("F18").Value = 69
$ws.Range("F25").Value = 28
$ws.Range("F28").Value = 77

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2609
$ws.Range("F6").Value = 1160
$ws.Range("F10").Value = 122
$ws.Range("F12").Value = 691

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2609
$ws.Range("F5").Value = 1160
$ws.Range("F8").Value = 122
$ws.Range("F9").Value = 682
$ws.Range("F10").Value = 691
$ws.Range("F11").Value = 610
$ws.Range("F14").Value = 38
$ws.Range("F15").Value = 6093
$ws.Range("F17").Value = 664
$ws.Range("F18").Value = 1066
$ws.Range("F19").Value = 27
$ws.Range("F20").Value = 264
$ws.Range("F23").Value = 601
$ws.Range("F28").Value = 69
$ws.Range("F29").Value = 1034
$ws.Range("F30").Value = 56
$ws.Range("F31").Value = 38
$ws.Range("F34").Value = 1363
$ws.Range("F36").Value = 28
$ws.Range("F39").Value = 1030
$ws.Range("F40").Value = 75
$ws.Range("F42").Value = 2095
$ws.Range("F44").Value = 204
$ws.Range("F45").Value = 41
$ws.Range("F48").Value = 3385
